$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places ("custom accuracy")
$row5 = @{
    "B5" = 10.09
    "C5" = 7.25
    "D5" = 0.89
    "E5" = 21.65
    "F5" = 17.94
    "G5" = 7.94
    "H5" = 34.01
    "I5" = 12.22
    "J5" = 5.35
    "K5" = 8.05
    "L5" = 8.67
    "M5" = 9.11
    "N5" = 2.54
    "O5" = 7.9
    "P5" = 11.19
    "Q5" = 6.76
    "R5" = 0.78
    "S5" = 0.52
    "T5" = 112.86
    "U5" = 22.19
    "V5" = 7.29
    "W5" = 14.78
    "X5" = 7.92
    "Y5" = 0.98
    "Z5" = 15.88
    "AA5" = 6.44
    "AB5" = 5.79
    "AC5" = 6.79
    "AD5" = 9.12
    "AE5" = 0.56
    "AF5" = 30.82
    "AG5" = 4.06
    "AH5" = 9.11
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove the last data row (row 6) entirely, shifting nothing below it up
$ws.Rows(6).Delete()
